# Update cryptos list values (price/volume) per the Jun 19 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.360.80"
$ws.Range("E2").Value = "  -0.66%  "
# Row 3
$ws.Range("D3").Value = "1.714.67"
$ws.Range("E3").Value = "  -1.34%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9972"
$ws.Range("E4").Value = "  -0.19%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.10"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9978"
$ws.Range("E6").Value = "  -0.17%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4864"
$ws.Range("E7").Value = "  -0.66%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2579"
$ws.Range("E8").Value = "  -3.32%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06165"
$ws.Range("E9").Value = "  -2.67%  "
# Row 10
$ws.Range("D10").Value = "1.717.04"
$ws.Range("E10").Value = "  -1.09%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06948"
$ws.Range("E11").Value = "  -1.37%  "
# Row 12
$ws.Range("E12").Value = "  -1.63%  "
# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5969"
$ws.Range("E13").Value = "  -2.30%  "
# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.467"
$ws.Range("E14").Value = "  -3.08%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.43"
$ws.Range("E15").Value = "  -1.33%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9976"
$ws.Range("E16").Value = "  -0.21%  "
# Row 17
$ws.Range("D17").Value = "26.351.85"
$ws.Range("E17").Value = "  -0.66%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9973"
$ws.Range("E18").Value = "  -0.21%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007093"
$ws.Range("E19").Value = "  -4.20%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.22"
$ws.Range("E20").Value = "  -2.61%  "
# Row 21
$ws.Range("D21").Value = "1.932.62"
$ws.Range("E21").Value = "  -1.03%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.406"
$ws.Range("E22").Value = "  -3.75%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.427"
$ws.Range("E23").Value = "  -3.36%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.040"
$ws.Range("E24").Value = "  -3.84%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.46"
$ws.Range("E25").Value = "  -3.17%  "
# Row 26
$ws.Range("E26").Value = "  -1.56%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.397"
$ws.Range("E27").Value = "  -0.82%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "105.74"
$ws.Range("E28").Value = "  -2.18%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.727"
$ws.Range("E29").Value = "  -2.27%  "
# Row 30
$ws.Range("E30").Value = "  -3.88%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07957"
$ws.Range("E31").Value = "  -1.19%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.605"
$ws.Range("E32").Value = "  -3.03%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04425"
$ws.Range("E33").Value = "  -3.38%  "
# Row 34
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9966"
$ws.Range("E34").Value = "  -0.23%  "
# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.600"
$ws.Range("E35").Value = "  -0.31%  "
# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9926"
$ws.Range("E36").Value = "  -1.70%  "
# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6172"
$ws.Range("E37").Value = "  -3.06%  "
# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9332"
$ws.Range("E38").Value = "  +4.25%  "
# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.973"
$ws.Range("E39").Value = "  -2.33%  "
# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.375"
$ws.Range("E40").Value = "  -1.02%  "
# Row 41
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9969"
$ws.Range("E41").Value = "  -0.75%  "
# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01473"
$ws.Range("E42").Value = "  -2.09%  "
# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.38"
$ws.Range("E43").Value = "  -3.14%  "
# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.422"
$ws.Range("E44").Value = "  +0.48%  "
# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3804"
$ws.Range("E45").Value = "  -2.46%  "
# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.810"
$ws.Range("E46").Value = "  -1.21%  "
# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1149"
$ws.Range("E47").Value = "  -3.08%  "
# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05348"
$ws.Range("E48").Value = "  -0.82%  "
# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.39"
$ws.Range("E49").Value = "  -0.40%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.681"
$ws.Range("E50").Value = "  -1.51%  "
# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.05"
$ws.Range("E51").Value = "  -1.46%  "
